# Revise config file handling
# Append three new daily log rows (144-146) to each of the four data sheets,
# mirroring the existing row layout (time, length, ID, actual-length, checksum
# plus their decimal counterparts).

$wb = $excel.ActiveWorkbook

# Values for column A (time, serial date numbers) shared by every sheet.
$dates = @(45930.43784722222, 45931.43920138889, 45932.43804398148)

# Per-sheet row data. Each entry is (B, C, D, E, F, G, H, I) for rows 144,145,146.
# Scientific-notation literals aren't accepted directly by the interpreter,
# so the big ID_DEC values are parsed from strings into doubles instead.
$idDecLFT1 = [double]"7.598631275147109e+23"
$idDecLFT2 = [double]"5.68432987514711e+23"
$idDecPLT1 = [double]"5.68631262647114e+23"
$idDecPLT2 = [double]"9.85046333984776e+23"

$sheetData = @{
    "DE_LFT_#1" = @(
        @("0x01,0x7c", "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,", "0x00,0xF8", "0x14", 380, $idDecLFT1, 248, 14),
        @("0x01,0x7c", "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,", "0x00,0xF4", "0x14", 380, $idDecLFT1, 244, 14),
        @("0x01,0x7c", "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,", "0x00,0xF0", "0x14", 380, $idDecLFT1, 240, 14)
    )
    "DE_LFT_#2" = @(
        @("0x01,0x7c", "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,", "0x00,0xF8", "0xe", 380, $idDecLFT2, 248, 14),
        @("0x01,0x7c", "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,", "0x00,0xF8", "0xe", 380, $idDecLFT2, 248, 14),
        @("0x01,0x7c", "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,", "0x00,0xF8", "0xe", 380, $idDecLFT2, 248, 14)
    )
    "DE_PLT_#1" = @(
        @("0x00,0x82", "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,", "0x00,0x6C", "0x7", 130, $idDecPLT1, 108, 7),
        @("0x00,0x82", "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,", "0x00,0x6C", "0x7", 130, $idDecPLT1, 108, 7),
        @("0x00,0x82", "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,", "0x00,0x6C", "0x7", 130, $idDecPLT1, 108, 7)
    )
    "DE_PLT_#2" = @(
        @("0x00,0x82", "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,", "0x00,0x6B", "0x3", 130, $idDecPLT2, 107, 3),
        @("0x00,0x82", "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,", "0x00,0x6B", "0x3", 130, $idDecPLT2, 107, 3),
        @("0x00,0x82", "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,", "0x00,0x6A", "0x3", 130, $idDecPLT2, 106, 3)
    )
}

foreach ($sheetName in $sheetData.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Re-use the date/time number format already applied to the preceding data rows.
    $dateFormat = $ws.Cells.Item(143, 1).NumberFormat

    $rows = $sheetData[$sheetName]
    for ($i = 0; $i -lt 3; $i++) {
        $r = 144 + $i
        $values = $rows[$i]

        $ws.Cells.Item($r, 1).Value = $dates[$i]
        $ws.Cells.Item($r, 1).NumberFormat = $dateFormat

        $ws.Cells.Item($r, 2).Value = $values[0]
        $ws.Cells.Item($r, 3).Value = $values[1]
        $ws.Cells.Item($r, 4).Value = $values[2]
        $ws.Cells.Item($r, 5).Value = $values[3]
        $ws.Cells.Item($r, 6).Value = $values[4]
        $ws.Cells.Item($r, 7).Value = $values[5]
        $ws.Cells.Item($r, 8).Value = $values[6]
        $ws.Cells.Item($r, 9).Value = $values[7]
    }
}
